# This document contains a table of three-digit-by-one-digit division
# "problems" (e.g. "337÷3=") that need updating to new values. Almost all
# of the old values are unique in the document, so a simple
# Find/Replace (ReplaceAll) is safe for those. One value, "289÷6=",
# appears twice in two different cells with two different replacements,
# so those two cells are targeted individually by scoping a Range to the
# specific table cell and doing a single (ReplaceOne) Find/Execute there.

$d = $word.ActiveDocument

function Replace-AllText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

function Replace-InCell($row, $col, $old, $new) {
    $cell = $d.Tables(1).Cell($row, $col)
    $cellRange = $cell.Range
    $scoped = $d.Range($cellRange.Start, $cellRange.End)
    $scoped.Find.Execute($old, $true, $false, $false, $false, $false, `
                          $true, 0, $false, $new, 1) | Out-Null
}

# Unique replacements (safe to do document-wide).
Replace-AllText "337÷3=" "160÷4="
Replace-AllText "806÷6=" "629÷8="
Replace-AllText "573÷7=" "991÷6="
Replace-AllText "782÷5=" "890÷8="
Replace-AllText "529÷3=" "491÷4="
Replace-AllText "659÷4=" "258÷4="
Replace-AllText "681÷8=" "334÷4="
Replace-AllText "842÷6=" "908÷5="
Replace-AllText "971÷2=" "652÷2="
Replace-AllText "335÷8=" "337÷2="
Replace-AllText "788÷8=" "483÷8="
Replace-AllText "828÷7=" "545÷3="
Replace-AllText "221÷8=" "969÷6="
Replace-AllText "406÷4=" "231÷8="
Replace-AllText "954÷7=" "378÷2="
Replace-AllText "285÷5=" "986÷2="
Replace-AllText "568÷3=" "880÷8="
Replace-AllText "502÷4=" "256÷3="
Replace-AllText "870÷7=" "404÷7="
Replace-AllText "623÷9=" "578÷2="
Replace-AllText "405÷9=" "330÷8="
Replace-AllText "234÷8=" "779÷8="
Replace-AllText "834÷2=" "731÷8="

# Duplicate "289÷6=" value -- handle each occurrence by its own cell.
Replace-InCell 9 1 "289÷6=" "528÷9="
Replace-InCell 17 5 "289÷6=" "298÷8="
